# Add "Closed Gt2 Belt 6mm * 110mm" line item to the BOM sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

$xlPasteFormats = -4122

# --- Formatting: clone row 7's look for the new row 13, then fix the K
#     column to match the bold style used by K11/K12 (even-row look). ---
$ws.Range("A7:K7").Copy()
$ws.Range("A13:K13").PasteSpecial($xlPasteFormats)
$ws.Range("K12").Copy()
$ws.Range("K13").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- New row 13 content ---
$ws.Range("A13").Value = "Closed Gt2 Belt 6mm * 110mm"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 9.99
$ws.Range("D13").Formula = "=C13*B13"
$ws.Range("E13").Value = "https://amzn.to/3NMXYm9"
$ws.Range("F13").Value = 15.79
$ws.Range("G13").Formula = "=F13*B13"
$ws.Range("H13").Value = "https://amzn.to/3CHSoLm"
$ws.Range("I13").Value = 0.81
$ws.Range("J13").Formula = "=I13"
$ws.Range("K13").Value = "https://s.click.aliexpress.com/e/_DEv5hNn"

# --- F12 previously blank, now explicitly zero ---
$ws.Range("F12").Value = 0

# --- Extend the three running totals down through the new row ---
$ws.Range("D15").Formula = "=SUM(D2:D13)"
$ws.Range("G15").Formula = "=SUM(G2:G13)"
$ws.Range("J15").Formula = "=SUM(J2:J13)"

# --- Selection moved to H28 when the edit was made ---
$ws.Range("H28").Select() | Out-Null

$wb.Save()
